# Fix contract addendum template:
# - Convert the three bulleted "Căn cứ..." paragraphs from list-numbered
#   bullets (numPr/numId=1) to literal "-" + tab text runs.
# - Normalize the "${base}" paragraph's indentation/tabs and drop its
#   leading tab run.
# - The now-unused numbering definition (numId=1) is cleared from every
#   paragraph that referenced it.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

$p1xml = '<w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:left="720" w:hanging="360"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>-</w:t><w:tab/><w:t>Căn cứ Bộ Luật Lao Động nước Cộng Hòa Xã Hội Chủ Nghĩa Việt Nam.</w:t></w:r></w:p>'
$p2xml = '<w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:left="720" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>-</w:t><w:tab/><w:t>Căn cứ nhu cầu lao động của Công ty TNHH DV Chấn Thanh</w:t></w:r></w:p>'
$p3xml = '<w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:left="720" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>-</w:t><w:tab/><w:t xml:space="preserve">Căn cứ hợp đồng lao động số </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>${number_contract}</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> ngày </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>${date_contract} giữa CÔNG TY TNHH DỊCH VỤ CHẤN THANH và Ông/Bà ${employee}</w:t></w:r></w:p>'
$p4xml = '<w:p><w:pPr><w:pStyle w:val="LOnormal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:left="360" w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:position w:val="0"/><w:sz w:val="20"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:highlight w:val="white"/><w:vertAlign w:val="baseline"/></w:rPr><w:t>${base}</w:t></w:r></w:p>'

$target1 = Find-ParagraphByText $d "Căn cứ Bộ Luật Lao Động"
$target1.Range.InsertXML($p1xml)

$target2 = Find-ParagraphByText $d "Căn cứ nhu cầu lao động"
$target2.Range.InsertXML($p2xml)

$target3 = Find-ParagraphByText $d "Căn cứ hợp đồng lao động số"
$target3.Range.InsertXML($p3xml)

$target4 = Find-ParagraphByText $d '${base}'
$target4.Range.InsertXML($p4xml)
